$d = $word.ActiveDocument

# 1) Fix spacing in fechaMemorandum placeholder: ", $      {incorporacion.fechaMemorandum}" -> ", ${incorporacion.fechaMemorandum}"
$d.Content.Find.Execute(", `$      {incorporacion.fechaMemorandum}", $true, $false, $false, $false, $false, $true, 1, $false, ", `${incorporacion.fechaMemorandum}", 2)

# 2) Swap the "a" / "destinatario" bookmark names (keep the same ranges/spans).
$bms = $d.Bookmarks
$bmA = $bms.Item("a")
$startA = $bmA.Start
$endA = $bmA.End
$bmDest = $bms.Item("destinatario")
$startD = $bmDest.Start
$endD = $bmDest.End

$bmDest.Delete()
$bmA.Delete()

$rangeD = $d.Range($startD, $endD)
$d.Bookmarks.Add("a", $rangeD)
$rangeA = $d.Range($startA, $endA)
$d.Bookmarks.Add("destinatario", $rangeA)

# 3) departamentoRef -> departamento, gerenciaRef -> gerencia
$d.Content.Find.Execute("`${puestoNuevo.departamentoRef}", $true, $false, $false, $false, $false, $true, 1, $false, "`${puestoNuevo.departamento}", 2)
$d.Content.Find.Execute("dependiente `${puestoNuevo.gerenciaRef} del Servicio", $true, $false, $false, $false, $false, $true, 1, $false, "dependiente `${puestoNuevo.gerencia} del Servicio", 2)

# 4) Merge "${" + "incorporacion" + ".gerenciaAbreviatura}" runs into "${puestoNuevo.gerenciaAbreviatura}"
#    by only rewriting the middle run's text; identically-formatted adjacent runs
#    coalesce automatically, while the preceding "    " run (same formatting) is left untouched.
$text = $d.Content.Text
$idx = $text.IndexOf("incorporacion.gerenciaAbreviatura")
$r = $d.Range($idx, $idx + "incorporacion".Length)
$r.Text = "puestoNuevo"

# 5) numeroTramite -> numeroHp
$d.Content.Find.Execute("`${incorporacion.numeroTramite}", $true, $false, $false, $false, $false, $true, 1, $false, "`${incorporacion.numeroHp}", 2)
